$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 428 (shifts the existing rows 428-527 down
# to 429-528, which matches how every cell's data "moved" by one row in the
# diff — the underlying per-row values are unchanged, only their row number
# increased by one).
$ws.Rows.Item(428).Insert()

# Populate the freshly inserted row 428 with its own data (same shape as the
# other rows in this "Apio" / Feria Lagunitas de Puerto Montt table).
$ws.Cells.Item(428, 1).Value2 = 4
$ws.Cells.Item(428, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(428, 3).Value2 = "Los Lagos"
$ws.Cells.Item(428, 4).Value2 = 45173
$ws.Cells.Item(428, 5).Value2 = 10
$ws.Cells.Item(428, 6).Value2 = 100112017
$ws.Cells.Item(428, 7).Value2 = "Apio"
$ws.Cells.Item(428, 8).Value2 = "Americana (o)"
$ws.Cells.Item(428, 9).Value2 = "Primera"
$ws.Cells.Item(428, 10).Value2 = 20
$ws.Cells.Item(428, 11).Value2 = 11000
$ws.Cells.Item(428, 12).Value2 = 11000
$ws.Cells.Item(428, 13).Value2 = 11000
$ws.Cells.Item(428, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(428, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(428, 16).Value2 = 1833
$ws.Cells.Item(428, 17).Value2 = 6
$ws.Cells.Item(428, 18).Value2 = "Hortaliza"
